# Auto-generated by analysis of OOXML diff for Carbuncle_Profits workbook.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) per leve row
# across all 8 job sheets, as produced by the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43: Growing Is Knowing
$ws.Range("H43").Value = 2324
$ws.Range("I43").Value = 2155
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 2155
$ws.Range("L43").Value = 3000
$ws.Range("M43").Value = -2086
$ws.Range("N43").Value = -3138

# Row 51: A Bile Business
$ws.Range("H51").Value = 3550
$ws.Range("I51").Value = 1160.2
$ws.Range("J51").Value = 4877.6665
$ws.Range("K51").Value = 1160.2
$ws.Range("L51").Value = 4877.6665
$ws.Range("M51").Value = -676.2
$ws.Range("N51").Value = -5845.6665

# Row 100: Asking for a Friend
$ws.Range("H100").Value = 27780556
$ws.Range("I100").Value = 55557676
$ws.Range("J100").Value = 3434.3333
$ws.Range("K100").Value = 55557676
$ws.Range("L100").Value = 3434.3333
$ws.Range("M100").Value = -55557135
$ws.Range("N100").Value = -4516.3333

# Row 113: Amaro Kart
$ws.Range("H113").Value = 8472.091
$ws.Range("I113").Value = 2851.6667
$ws.Range("J113").Value = 10579.75
$ws.Range("K113").Value = 2851.6667
$ws.Range("L113").Value = 10579.75
$ws.Range("M113").Value = 402.3332999999998
$ws.Range("N113").Value = -17087.75

# Row 117: A Greater Grimoire
$ws.Range("H117").Value = 34980
$ws.Range("J117").Value = 34980
$ws.Range("L117").Value = 34980
$ws.Range("N117").Value = -44158

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1588.2667
$ws.Range("I137").Value = 1206.4166
$ws.Range("J137").Value = 3115.6667
$ws.Range("K137").Value = 3619.2498
$ws.Range("L137").Value = 9347.000100000001
$ws.Range("M137").Value = -1069.2498
$ws.Range("N137").Value = -14447.0001

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 17250
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 24: A Firm Hand
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

# Row 32: Ingot We Trust
$ws.Range("H32").Value = 3322.7532
$ws.Range("I32").Value = 2413.2297
$ws.Range("K32").Value = 2413.2297
$ws.Range("M32").Value = -2126.2297

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 1218.5264
$ws.Range("I45").Value = 810.375
$ws.Range("J45").Value = 1515.3636
$ws.Range("K45").Value = 810.375
$ws.Range("L45").Value = 1515.3636
$ws.Range("M45").Value = -433.375
$ws.Range("N45").Value = -2269.3636

# Row 96: The Gauntlet Is Cast
$ws.Range("H96").Value = 18958
$ws.Range("J96").Value = 18958
$ws.Range("L96").Value = 18958
$ws.Range("N96").Value = -24450

# Row 97: Ore for Me
$ws.Range("H97").Value = 1600.2759
$ws.Range("I97").Value = 1251.5
$ws.Range("K97").Value = 1251.5
$ws.Range("M97").Value = -755.5

# Row 100: En Garde and on Guard
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 80: Unbreaker
$ws.Range("H80").Value = 4227.278
$ws.Range("I80").Value = 809.1
$ws.Range("J80").Value = 8500
$ws.Range("K80").Value = 809.1
$ws.Range("L80").Value = 8500
$ws.Range("M80").Value = 188.9
$ws.Range("N80").Value = -10496

# Row 83: Attack on Titanium (L)
$ws.Range("H83").Value = 4227.278
$ws.Range("I83").Value = 809.1
$ws.Range("J83").Value = 8500
$ws.Range("K83").Value = 4045.5
$ws.Range("L83").Value = 42500
$ws.Range("M83").Value = 946.5
$ws.Range("N83").Value = -52484

# Row 107: The Gold Experience
$ws.Range("H107").Value = 2214.1052
$ws.Range("I107").Value = 1837.3334
$ws.Range("K107").Value = 1837.3334
$ws.Range("M107").Value = 82.66660000000002

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2568.7917
$ws.Range("I134").Value = 1431.2142
$ws.Range("J134").Value = 4161.4
$ws.Range("K134").Value = 4293.642599999999
$ws.Range("L134").Value = 12484.2
$ws.Range("M134").Value = -1758.642599999999
$ws.Range("N134").Value = -17554.2

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 1311.069
$ws.Range("I16").Value = 1335.1305
$ws.Range("J16").Value = 1218.8334
$ws.Range("K16").Value = 1335.1305
$ws.Range("L16").Value = 1218.8334
$ws.Range("M16").Value = -1048.1305
$ws.Range("N16").Value = -1792.8334

# Row 28: Militia on My Mind
$ws.Range("H28").Value = 49385.8
$ws.Range("J28").Value = 49385.8
$ws.Range("L28").Value = 49385.8
$ws.Range("N28").Value = -49875.8

# Row 31: Wall Not Found
$ws.Range("H31").Value = 2370.3333
$ws.Range("I31").Value = 1887.0741
$ws.Range("J31").Value = 3095.2222
$ws.Range("K31").Value = 1887.0741
$ws.Range("L31").Value = 3095.2222
$ws.Range("M31").Value = -1592.0741
$ws.Range("N31").Value = -3685.2222

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2370.3333
$ws.Range("I34").Value = 1887.0741
$ws.Range("J34").Value = 3095.2222
$ws.Range("K34").Value = 1887.0741
$ws.Range("L34").Value = 3095.2222
$ws.Range("M34").Value = -1685.0741
$ws.Range("N34").Value = -3499.2222

# Row 86: Birch, Please
$ws.Range("H86").Value = 2529.8235
$ws.Range("I86").Value = 2681.9167
$ws.Range("J86").Value = 2164.8
$ws.Range("K86").Value = 2681.9167
$ws.Range("L86").Value = 2164.8
$ws.Range("M86").Value = -1558.9167
$ws.Range("N86").Value = -4410.8

# Row 89: Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 2529.8235
$ws.Range("I89").Value = 2681.9167
$ws.Range("J89").Value = 2164.8
$ws.Range("K89").Value = 13409.5835
$ws.Range("L89").Value = 10824
$ws.Range("M89").Value = -7793.583500000001
$ws.Range("N89").Value = -22056

# Row 113: Patient Patients
$ws.Range("H113").Value = 1311.069
$ws.Range("I113").Value = 1335.1305
$ws.Range("J113").Value = 1218.8334
$ws.Range("K113").Value = 1335.1305
$ws.Range("L113").Value = 1218.8334
$ws.Range("M113").Value = 834.8695
$ws.Range("N113").Value = -5558.8334

# Row 114: Ground to a Halt
$ws.Range("H114").Value = 49456
$ws.Range("J114").Value = 49456
$ws.Range("L114").Value = 49456
$ws.Range("N114").Value = -58134

# Row 123: A Real Grind
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2152.8708
$ws.Range("I132").Value = 1124.1818
$ws.Range("J132").Value = 4667.4443
$ws.Range("K132").Value = 3372.5454
$ws.Range("L132").Value = 14002.3329
$ws.Range("M132").Value = -842.5454
$ws.Range("N132").Value = -19062.3329

$ws = $wb.Worksheets.Item("CUL")
# Row 82: Persuasion of a Higher Power
$ws.Range("H82").Value = 170899.83
$ws.Range("J82").Value = 204879.8
$ws.Range("L82").Value = 614639.3999999999
$ws.Range("N82").Value = -615451.3999999999

# Row 85: Loaves and Fishes (L)
$ws.Range("H85").Value = 170899.83
$ws.Range("J85").Value = 204879.8
$ws.Range("L85").Value = 614639.3999999999
$ws.Range("N85").Value = -617447.3999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 2270
$ws.Range("I80").Value = 1566.6666
$ws.Range("J80").Value = 2571.4285
$ws.Range("K80").Value = 1566.6666
$ws.Range("L80").Value = 2571.4285
$ws.Range("M80").Value = -568.6666
$ws.Range("N80").Value = -4567.4285

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 2270
$ws.Range("I83").Value = 1566.6666
$ws.Range("J83").Value = 2571.4285
$ws.Range("K83").Value = 7833.333000000001
$ws.Range("L83").Value = 12857.1425
$ws.Range("M83").Value = -2841.333000000001
$ws.Range("N83").Value = -22841.1425

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 1986.0209
$ws.Range("I122").Value = 1767.921
$ws.Range("K122").Value = 5303.763
$ws.Range("M122").Value = -2853.763

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 1966.2858
$ws.Range("I126").Value = 1620.6957
$ws.Range("J126").Value = 2628.6667
$ws.Range("K126").Value = 4862.0871
$ws.Range("L126").Value = 7886.000100000001
$ws.Range("M126").Value = -2392.0871
$ws.Range("N126").Value = -12826.0001

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 69559.87
$ws.Range("I7").Value = 113288.11
$ws.Range("J7").Value = 3967.5
$ws.Range("K7").Value = 113288.11
$ws.Range("L7").Value = 3967.5
$ws.Range("M7").Value = -113176.11
$ws.Range("N7").Value = -4191.5

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 1544.1428
$ws.Range("I61").Value = 1400
$ws.Range("J61").Value = 1601.8
$ws.Range("K61").Value = 1400
$ws.Range("L61").Value = 1601.8
$ws.Range("M61").Value = -1198
$ws.Range("N61").Value = -2005.8

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 4164.6206
$ws.Range("I93").Value = 4773.5835
$ws.Range("J93").Value = 1241.6
$ws.Range("K93").Value = 4773.5835
$ws.Range("L93").Value = 1241.6
$ws.Range("M93").Value = -3525.5835
$ws.Range("N93").Value = -3737.6

# Row 113: Peace in Rest
$ws.Range("H113").Value = 1544.1428
$ws.Range("I113").Value = 1400
$ws.Range("J113").Value = 1601.8
$ws.Range("K113").Value = 1400
$ws.Range("L113").Value = 1601.8
$ws.Range("M113").Value = 770
$ws.Range("N113").Value = -5941.8

# Row 126: Battered Books
$ws.Range("H126").Value = 69559.87
$ws.Range("I126").Value = 113288.11
$ws.Range("J126").Value = 3967.5
$ws.Range("K126").Value = 339864.33
$ws.Range("L126").Value = 11902.5
$ws.Range("M126").Value = -337394.33
$ws.Range("N126").Value = -16842.5

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display
$ws.Range("H96").Value = 1434.875
$ws.Range("I96").Value = 1413.1666
$ws.Range("K96").Value = 1413.1666
$ws.Range("M96").Value = -40.16660000000002

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 34474.8
$ws.Range("I126").Value = 46324.5
$ws.Range("J126").Value = 1888.125
$ws.Range("K126").Value = 138973.5
$ws.Range("L126").Value = 5664.375
$ws.Range("M126").Value = -136503.5
$ws.Range("N126").Value = -10604.375
